$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 6 cells whose text content actually changed ---

$ws.Range("C44").Value = '"operacion":{"vendedor":{"cuit":"20000001695","cbu":"0000489700000000000000","banco":"000","recurrencia":true,"prestacion":"Prestacion10"},"comprador":{"cuit":"20000001695","cuenta":{"cbu":"0000489700000000000017"}}}}|"operacion":{"comprador":{"cuit":"20000001695","cuenta":{"cbu":"0000489700000000000017"}}}}|"operacion":{"comprador":{"cuit":"20000001695","cuenta":{"cbu":"9980000400000000000758"}}}}|"id":"debin.id","aviso":"all","producto":"responder"'

$ws.Range("C49").Value = '"operacion":{"vendedor":{"cuit":"20000001652","cbu":"0000494100000000000000","banco":"000","recurrencia":true,"prestacion":"Prestacion6"},"comprador":{"cuit":"20000001652","cuenta":{"cbu":"0000494100000000000017"}}}}|"operacion":{"comprador":{"cuit":"20000001652","cuenta":{"cbu":"0000494100000000000017"}}}}|"operacion":{"comprador":{"cuit":"20000001652","cuenta":{"cbu":"9980000400000000000673"}}}}|"id":"debin.id","aviso":"all","producto":"responder"'

$ws.Range("C57").Value = '"operacion":{"vendedor":{"cuit":"20000001725","cbu":"0000484200000000000000","banco":"000","sucursal":"0484"},"comprador":{"cuenta":{"cbu":"0000484200000000000017"},"cuit":"20000001725"},"detalle":{"importe":10,"id_billetera":484}}}|"operacion":{"comprador":{"cuit":"20000001725","cuenta":{"cbu":"9980000400000000000802"}},"detalle":{"importe":10}}}|"id":"debin.id","aviso":"all","producto":"responder"'

$ws.Range("C61").Value = '"operacion":{"vendedor":{"cuit":"20000001776","cbu":"0000476700000000000000","banco":"000","sucursal":"0476"},"comprador":{"cuenta":{"cbu":"0000476700000000000017"},"cuit":"20000001776"},"detalle":{"importe":1000,"id_billetera":476}}}|"operacion":{"comprador":{"cuit":"20000001776","cuenta":{"cbu":"9980000400000000000895"}},"detalle":{"importe":1000}}}|"id":"debin.id","aviso":"all","producto":"responder"'

$ws.Range("C64").Value = '"operacion":{"vendedor":{"cuit":"20000001784","cbu":"0000474300000000000000","banco":"000","sucursal":"0474"},"comprador":{"cuenta":{"cbu":"0000474300000000000017"},"cuit":"20000001784"},"detalle":{"importe":1000,"id_billetera":474}}}|"operacion":{"comprador":{"cuit":"20000001784","cuenta":{"cbu":"9980000400000000000918"}},"detalle":{"importe":1000}}}|"id":"debin.id","aviso":"all","producto":"responder"'

$ws.Range("C66").Value = '"operacion":{"vendedor":{"cuit":"20000001725","cbu":"0000484200000000000000","banco":"000","sucursal":"0484"},"comprador":{"cuenta":{"cbu":"0000484200000000000017"},"cuit":"20000001725"},"detalle":{"importe":1000,"id_billetera":484}}}|"operacion":{"comprador":{"cuit":"20000001725","cuenta":{"cbu":"9980000400000000000802"}},"detalle":{"importe":1000}}}|"operacion_original":{"detalle":{"importe":1000},"vendedor":{"cuit":"20000001725","cbu":"0000484200000000000000"}}}|"id":"debin.id","aviso":"all","producto":"responder"'

# --- Row 49 shrinks from 90 to 75 once the shorter text reflows ---
$ws.Rows.Item(49).RowHeight = 75

# --- View state: scroll position + selection moved down the sheet ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A61").Select() | Out-Null
